$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 500
$ws.Range("C2").Value = 508
$ws.Range("D2").Value = 669
$ws.Range("E2").Value = 668
$ws.Range("F2").Value = 829
$ws.Range("G2").Value = 838
$ws.Range("B3").Value = 942
$ws.Range("C3").Value = 948
$ws.Range("D3").Value = 1054
$ws.Range("E3").Value = 1055
$ws.Range("F3").Value = 1158
$ws.Range("G3").Value = 1163
$ws.Range("C4").Value = 20.69
$ws.Range("F4").Value = 35.97
$ws.Range("C5").Value = 20.69
$ws.Range("F5").Value = 35.96
$ws.Range("B6").Value = 2.79
$ws.Range("G6").Value = 3.79
$ws.Range("B7").Value = 2.8
$ws.Range("G7").Value = 3.83
$ws.Range("G8").Value = 6.33
$ws.Range("G9").Value = 6.27
$ws.Range("B10").Value = 4000
$ws.Range("C10").Value = 4398.75
$ws.Range("D10").Value = 7976.39
$ws.Range("E10").Value = 7898.23
$ws.Range("F10").Value = 12079.54
$ws.Range("G10").Value = 12999.99
$ws.Range("C11").Value = 4399.27
$ws.Range("D11").Value = 7978.32
$ws.Range("E11").Value = 7900.75
$ws.Range("F11").Value = 12082.2
$ws.Range("C16").Value = 101
$ws.Range("D16").Value = 120
$ws.Range("D17").Value = 119.99
$ws.Range("F17").Value = 139
$ws.Range("F18").Value = 0.35
$ws.Range("C19").Value = 0.15
$ws.Range("F19").Value = 0.35
$ws.Range("B22").Value = 7983.03
$ws.Range("C22").Value = 9307.61
$ws.Range("D22").Value = 11632.73
$ws.Range("E22").Value = 11580.81
$ws.Range("F22").Value = 14310.4
$ws.Range("G22").Value = 16765.13
$ws.Range("B23").Value = 9761.47
$ws.Range("C23").Value = 12837.3
$ws.Range("D23").Value = 16568.43
$ws.Range("E23").Value = 16461.64
$ws.Range("F23").Value = 20920.88
$ws.Range("G23").Value = 27570.41
$ws.Range("B24").Value = 0.88
$ws.Range("G24").Value = 6.91
